# Scene 76 revision: Prim and the narrator now head towards the exit
# "silently" after she helps him to his feet.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "the two of us head towards the exit",  # Find what
    $true,                                   # MatchCase
    $false,                                  # MatchWholeWord
    $false,                                  # MatchWildcards
    $false,                                  # MatchSoundsLike
    $false,                                  # MatchAllWordForms
    $true,                                   # Forward
    1,                                        # Wrap (wdFindContinue)
    $false,                                  # Format
    "the two of us silently head towards the exit",  # Replace with
    2                                         # Replace (wdReplaceAll)
)
